$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New column D: width, used later by D11/D12/D18 (wrap-text style emerges
# automatically from WrapText usage below - style index 5 gets created).
# ColumnWidth of 84.16666666666667 serializes to the OOXML width="85"
# (stored = (round(chars*6)+5)/6).
# ---------------------------------------------------------------------------
$ws.Columns(4).ColumnWidth = 84.16666666666667

# ---------------------------------------------------------------------------
# New rows 37-41: two people highlighted in red (fill 255 = RGB(255,0,0),
# matching existing style s="1"), plus three plain names.
# Values are written in this specific order so that newly created shared
# strings land at the same table indices as in the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("A38").Value = "茅勤"
$ws.Range("B38").Interior.Color = 255

$ws.Range("A37").Value = "朱保全"
$ws.Range("B37").Interior.Color = 255

$ws.Range("A39").Value = "曹京明"
$ws.Range("A40").Value = "吴轶秦"
$ws.Range("A41").Value = "郭散皞"

# ---------------------------------------------------------------------------
# Row 12 gets a tall row and a new wrapped note in D12.
# ---------------------------------------------------------------------------
$ws.Range("D12").Value = "领导，之前跟您汇报过一次，我开发了一个辅助车位定价的系统，现在可以运转。"
$ws.Range("D12").WrapText = $true
$ws.Rows(12).RowHeight = 135.75

# ---------------------------------------------------------------------------
# Row 18 gets a taller row and a new wrapped note in D18.
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = "我是北京区域设计部的李谦，曾经跟您见过几面。我最近研发了一个服务于车位定价的系统。这个系统对咱们现在的工作可能有一些价值。冒昧给您发一些资料，争取一些您的关注。这个系统采用了一些技术手段，可以给项目提供精细的车位定价方案。可以作为管理抓手、提高定价方案的质量、节约管理人力、提高透明度；解决车位销售中管理难度大、价值流失和透明度不足导致的客户信任问题。推荐给您是我考虑这个系统可能需要比较大的体量以产生足够价值。我给您附上一个介绍，期待您的指导和建议。"
$ws.Range("D18").WrapText = $true
$ws.Rows(18).RowHeight = 90

# ---------------------------------------------------------------------------
# D11 already has text (代销公司联系人) - just give it the wrap-text style.
# ---------------------------------------------------------------------------
$ws.Range("D11").WrapText = $true

# ---------------------------------------------------------------------------
# Update the view: select G12 (this also drops the old topLeftCell scroll
# position that pinned the view down near A52).
# ---------------------------------------------------------------------------
$ws.Range("G12").Select()
